$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scrape row appended below the existing two data rows.
$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("B4").Value = "IXIC"
$ws.Range("C4").Value = "NASDAQ Composite"
$ws.Range("D4").Value = "Trading in Progress"

# Current Price / Stock Change are stored as text (not numbers) in the
# source data, so force text entry and then strip the resulting style
# override back to Normal so no extra formatting is introduced.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "12714.06"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "+71.05  +0.56%"
$ws.Range("F4").Style = "Normal"

$ws.Range("G4").Value = 12842.96
$ws.Range("H4").Value = 12750.47
$ws.Range("I4").Value = 14446.55
$ws.Range("J4").Value = 0.91
$ws.Range("K4").Value = 12691.56
$ws.Range("L4").Value = 12643.01
$ws.Range("M4").Value = 10207.47
$ws.Range("N4").Value = 0.012
$ws.Range("O4").Value = 1873000000
